$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$ws.Range("L2").Value = 1.09
$ws.Range("L3").Value = 1.13
$ws.Range("L4").Value = 1.14
$ws.Range("L5").Value = 1.12
$ws.Range("L6").Value = 0.84
$ws.Range("L7").Value = 0.99
$ws.Range("L8").Value = 1.19
$ws.Range("L9").Value = 1.13
$ws.Range("L10").Value = 1.18
$ws.Range("L11").Value = 1.18
$ws.Range("L12").Value = 1.2
$ws.Range("L13").Value = 0.96
$ws.Range("L14").Value = 1.14
$ws.Range("L15").Value = 1.15
$ws.Range("L16").Value = 0.9399999999999999
$ws.Range("L17").Value = 0.9
